$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text/string storage for the Price column (D) so numeric-looking
# strings like "1.002" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '31.267.73'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '1.966.00'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '246.07'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '0.4919'
$ws.Range('E7').Value = '  +2.33%  '
$ws.Range('D8').Value = '44.87'
$ws.Range('E8').Value = '  -0.73%  '
$ws.Range('D9').Value = '0.2987'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').Value = '0.06891'
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('D11').Value = '19.40'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('D12').Value = '108.38'
$ws.Range('E12').Value = '  -4.11%  '
$ws.Range('D13').Value = '1.943.33'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').Value = '0.07769'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').Value = '5.481'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('D16').Value = '0.7175'
$ws.Range('E16').Value = '  +3.73%  '
$ws.Range('D17').Value = '287.29'
$ws.Range('E17').Value = '  -3.60%  '
$ws.Range('D18').Value = '31.168.65'
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('D19').Value = '0.000007814'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').Value = '13.29'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '2.198.68'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = '5.534'
$ws.Range('E23').Value = '  -2.56%  '
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '6.595'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = '9.877'
$ws.Range('E26').Value = '  +1.35%  '
$ws.Range('D27').Value = '169.78'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = '20.36'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').Value = '2.219'
$ws.Range('E29').Value = '  +1.87%  '
$ws.Range('D30').Value = '0.1057'
$ws.Range('E30').Value = '  -2.84%  '
$ws.Range('D31').Value = '1.436'
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').Value = '1.592'
$ws.Range('D33').Value = '4.634'
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('D34').Value = '4.467'
$ws.Range('E34').Value = '  +0.74%  '
$ws.Range('D35').Value = '0.04993'
$ws.Range('E35').Value = '  -1.34%  '
$ws.Range('D36').Value = '0.7638'
$ws.Range('E36').Value = '  -2.31%  '
$ws.Range('D37').Value = '1.184'
$ws.Range('E37').Value = '  +1.65%  '
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').Value = '0.02045'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').Value = '2.712'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').Value = '2.197'
$ws.Range('E41').Value = '  +7.66%  '
$ws.Range('D42').Value = '6.442'
$ws.Range('E42').Value = '  +7.90%  '
$ws.Range('D43').Value = '0.4559'
$ws.Range('E43').Value = '  +1.86%  '
$ws.Range('D44').Value = '109.94'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').Value = '0.8844'
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('D46').Value = '72.48'
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('D47').Value = '8.056'
$ws.Range('E47').Value = '  +8.62%  '
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').Value = '9.440'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.1273'
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '0.2643'
$ws.Range('E51').Value = '  +3.11%  '

# Restore the default (unstyled) cell style now that values are set,
# so no stray number-format style lingers on the Price column.
$ws.Range("D2:D51").Style = "Normal"

